# Review tracker update: append a new review-comment row (row 9) to Sheet1
# and move the reviewer's selection/focus down to where the new entry was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (S.No=8) mirroring the existing "Open" review rows (7 & 8):
#   Script Name        = All
#   Package/Location   = src/test/java/com/comcast/century/cm/pages (two lines)
#   Review Comment     = new comment about unhandled link/button clicks
#   Assigned           = Rijin/Kesavan/Jatin
#   Status             = Harsh
#   Date reported      = Open
#   Date fixed         = 9/30/2016 (blank, matching row 8's pattern)
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "All"
$ws.Cells.Item(9, 3).Value = "src/test/java/com/comcast/century/cm/pages`nsrc/test/java/com/comcast/century/cm/pages"
$ws.Cells.Item(9, 4).Value = "I still see that there are few link or button clicks that are not handled with page.java metho (Iclick)  `nfor.e.g. - in AccountTabPageCM.java. Line num 156, line 189"
$ws.Cells.Item(9, 5).Value = "Rijin/Kesavan/Jatin"
$ws.Cells.Item(9, 6).Value = "Harsh"
$ws.Cells.Item(9, 7).Value = "Open"

# Date fixed column: store the same serial date Excel would (30-Sep-2016),
# reusing row 8's date-cell number format so it renders identically (d-mmm style).
$ws.Cells.Item(9, 8).Value = 42643
$ws.Cells.Item(9, 8).NumberFormat = $ws.Cells.Item(8, 8).NumberFormat

# Match the row height Excel computed for this wrapped, multi-line comment.
$ws.Rows.Item(9).RowHeight = 75

# Reviewer's cursor ends up on I9 after entering the row.
$ws.Range("I9").Select() | Out-Null
